# RegNetAgents conference poster: revise timing language for defensibility.
# Replaces unsupported quantitative timing claims with qualitative language,
# touching three text boxes on slide 1 (abstract/key-innovation, performance
# impact bullets, and the challenge bullet list).

$arrow  = [char]0x2192   # "->" glyph used between "Hours"/"Manual" and the result
$bullet = [char]0x2022   # "*" separator glyph used in the KEY INNOVATION line

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 4 "TextBox 4": ABSTRACT body -----------------------------------
$sh4 = $s.Shapes.Item(4)
$tr4 = $sh4.TextFrame.TextRange
$para4 = $tr4.Paragraphs(2, 1)

$para4.Runs(1, 1).Text = "RegNetAgents is a validated LLM-powered multi-agent AI framework that automates gene regulatory network analysis through intelligent workflow orchestration, transforming labor-intensive manual processes into second-scale automated analysis. The system deploys four specialized LLM-powered domain agents (cancer biology, drug discovery, clinical relevance, systems biology) using local language model inference (Ollama/llama3.1:8b) to generate scientific insights with rationales, with graceful fallback to rule-based heuristics for reliability. Built on pre-computed regulatory networks derived from 500K+ single-cell RNA-seq profiles from CellxGene Data Portal (processed via ARACNe algorithm), perturbation analysis ranks therapeutic targets using network centrality metrics (PageRank, out-degree centrality), analyzing all upstream regulators comprehensively and successfully identifying experimentally validated regulators. Framework validation on colorectal cancer biomarkers showed 100% concordance with published literature across five genes and complete perturbation analysis of 99 regulators. Perturbation analysis successfully identified experimentally validated TP53 regulators (WWTR1, YAP1, CHD4 from Hippo pathway) alongside novel testable hypotheses (RBPMS, PRRX2), demonstrating reliable hypothesis generation for experimental prioritization. Natural language interface via Claude Desktop makes sophisticated gene analysis accessible without programming."

$para4.Runs(2, 1).Text = "KEY INNOVATION: LLM-Powered Agents with Scientific Rationales $bullet Local Inference (Ollama) $bullet Manual $arrow Seconds (15-62 sec with LLM, <1 sec rule-based) $bullet 4 Parallel Domain Agents $bullet Complete Perturbation Analysis (All Regulators) $bullet Conversational Interface $bullet Graceful Fallback Architecture"

# --- Shape 6 "TextBox 6": WHY THIS IS NOVEL / Performance Impact ---------
$sh6 = $s.Shapes.Item(6)
$tr6 = $sh6.TextFrame.TextRange

$tr6.Paragraphs(12, 1).Runs(1, 1).Text = "Single gene (rule-based): Manual $arrow 0.68 sec"
$tr6.Paragraphs(13, 1).Runs(1, 1).Text = "Single gene (LLM-powered): Manual $arrow 15 sec"
$tr6.Paragraphs(14, 1).Runs(1, 1).Text = "5 genes (LLM-powered): Manual $arrow 62 sec"

# --- Shape 8 "TextBox 8": THE CHALLENGE -----------------------------------
$sh8 = $s.Shapes.Item(8)
$tr8 = $sh8.TextFrame.TextRange

$tr8.Paragraphs(3, 1).Runs(1, 1).Text = "Manual analysis across multiple domains (cancer, drug, clinical) requires labor-intensive effort per gene"
